# Insert a new weekly price record as row 100 on the only worksheet.
# Inserting a whole row shifts the existing rows 100-122 down to 101-123
# (carrying their formatting/styles with them), matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("100:100").Insert()

$ws.Range("A100").Value = 6
$ws.Range("B100").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C100").Value = "Metropolitana"
$ws.Range("D100").Value = 44511
$ws.Range("E100").Value = 13
$ws.Range("F100").Value = 100112001
$ws.Range("G100").Value = "Berenjena"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 270
$ws.Range("K100").Value = 7000
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 7556
$ws.Range("N100").Value = "`$/caja 50 unidades"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 151
$ws.Range("Q100").Value = 50
$ws.Range("R100").Value = "Hortaliza"
